# BackupCodes.xlsx - CRDC backup codes update
#
# Two already-consumed backup codes (B5M3N2BW0MA3, XA43JHAGQ8V3) are
# removed from the top of the list and the now-unused blank slots at
# A15/A16 are cleared, so the remaining/added codes compact up into
# A2:A4 while the trailing codes (A17:A21) keep their original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Codes that used to live at A4/A15/A16 move up to fill the gap left by
# the two removed (used) codes that were at A2/A3.
$ws.Range("A2").Value = "TK50GFCXFHCN"
$ws.Range("A3").Value = "9M1A883VTX21"
$ws.Range("A4").Value = "D2ANGGG71FGC"

# The rows these values used to occupy are now empty.
$ws.Range("A15:A16").ClearContents()

# A17:A21 keep their existing codes (8AW7QA18SBTA, EA5XZ049QR7S,
# PTV3TSFPBF6W, NHYK5008HQDA, ZADGNDVPP03M) untouched.

# Matches the saved selection state in the workbook after the edit.
[void]$ws.Range("A6").Select()
